$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date serial values for rows 2-5 advance by one day
# (2023-09-14 -> 2023-09-15), i.e. 45183 -> 45184.
$ws.Range("C2").Value = 45184
$ws.Range("C3").Value = 45184
$ws.Range("C4").Value = 45184
$ws.Range("C5").Value = 45184
